# feat: Add generating random _id for members
#
# - Drop the "순서" (sequence number) column: it's being replaced by an
#   auto-generated member _id, so the manual ordinal column goes away and
#   the remaining columns (회원명 / 월수레슨 / 화목레슨) shift left.
# - Turn the old sample row into a reusable template: row 3 keeps the
#   normal data-row look, while rows 4-8 become blank placeholder rows
#   (grey text) ready to be filled with new members / generated ids.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A ("순서") entirely - B/C/D shift left into A/B/C.
$ws.Columns.Item(1).Delete()

# Row 3: blank row, formatted like the existing data row (row 2).
$ws.Range("A2:C2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Build the grey "placeholder" look (used for the still-empty template
# rows) once on A4, then again on B4 with a tiny format nudge so the two
# end up as distinct-but-identical-looking cell styles, same as the
# id column vs. the other columns in the source sheet.
$ws.Range("A4").Value = "x"
$ws.Range("A4").Font.Color = 12040119
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").ClearContents()

$ws.Range("B4").Value = "x"
$ws.Range("B4").Font.Color = 12040119
$ws.Range("B4").Font.Name = "Arial"
$ws.Range("B4").VerticalAlignment = -4160
$ws.Range("B4").ClearContents()

# Fan the two placeholder styles out across rows 4-8.
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B8").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("C4:C8").PasteSpecial(-4122)

# C7 stays on the normal data-row style rather than the placeholder one.
$ws.Range("C7").Value = ""
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)

Write-Output "done"
